{"js": "// Update progress marker: \"In progress\" (blue, accent1) -> \"Done\" (green 00B050).\nconst body = context.document.body;\n\nconst statusResults = body.search(\"In progress\", { matchCase: true, matchWholeWord: false });\nstatusResults.load(\"items\");\nawait context.sync();\n\nif (statusResults.items.length > 0) {\n  const statusRange = statusResults.items[0];\n  statusRange.font.color = \"#00B050\";\n  statusRange.insertText(\"Done\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// The \"_GoBack\" bookmark (Word's last-edit-position marker) moves from the\n// end of the status run to the point in the \"planets placed randomly\"\n// bullet where the most recent typing happened (right after \"...have pl\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst anchorResults = body.search(\"It would be cool to have pl\", { matchCase: true, matchWholeWord: false });\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n  const anchorRange = anchorResults.items[0];\n  const afterAnchor = anchorRange.getRange(Word.RangeLocation.after);\n  afterAnchor.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Update progress marker: \"In progress\" (blue, accent1) -> \"Done\" (green 00B050).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"In progress\"\n$find.MatchCase = $true\n$found = $find.Execute()\nif ($found) {\n    $range.Text = \"Done\"\n    $range.Font.Bold = $true\n    $range.Font.Color = 5287936   # RGB(0,176,80) => 00B050, stored as BGR\n}\n\n# The \"_GoBack\" bookmark (Word's last-edit-position marker) moves from the\n# end of the status run to the point in the \"planets placed randomly\"\n# bullet where the most recent typing happened (right after \"...have pl\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$anchorRange = $d.Content\n$anchorFind = $anchorRange.Find\n$anchorFind.Text = \"It would be cool to have pl\"\n$anchorFind.MatchCase = $true\n$anchorFound = $anchorFind.Execute()\nif ($anchorFound) {\n    $bmRange = $d.Range($anchorRange.End, $anchorRange.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
